$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $ok = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    return $ok
}

# ---------------------------------------------------------------------------
# Change 1: collapse the "letra tamaño" run split (drops the gramStart/gramEnd
# proofErr markers by merging the three runs " (", "letra tamaño",
# " 10, se sugiere Arial o " into a single run).
# ---------------------------------------------------------------------------
$old1 = " (letra tamaño 10, se sugiere Arial o "
$new1 = " (letra tamaño 10, se sugiere Arial o "
Replace-Text $old1 $new1 | Out-Null

# ---------------------------------------------------------------------------
# Change 2: append the new sentence about the heart-rate / temperature
# sensors right after "... procesamiento y envió de datos al dispositivo
# externo" in the "Objetivo General" paragraph.
# ---------------------------------------------------------------------------
$old2 = " procesamiento y envió de datos al dispositivo externo"
$new2 = " procesamiento y envió de datos al dispositivo externo" + `
    ". Para mediar la frecuencia cardiaca se " + `
    "usará" + `
    " un sensor de ritmo cardiaco el cual combina un sensor de pulsos ópticos con un circuito de amplificación y cancelación de ruido lo que permite obtener lecturas confiables de las pulsaciones del corazón, para la medición de la temperatura nos aprovecharemos de los materiales de los sensores" + `
    "," + `
    " los cuales varían su resistencia " + `
    "eléctrica por ende podremos detecta varianzas en los cambios de temperatura" + `
    " y registrando la temperatura obtenida por el sensor" + `
    "."
Replace-Text $old2 $new2 | Out-Null

# ---------------------------------------------------------------------------
# Change 3: merge the four runs of the "La interpretación..." paragraph into
# a single run, then delete all the empty trailing paragraphs that followed
# it (down to the sectPr).
# ---------------------------------------------------------------------------
$old3 = "La interpretación de los parámetros permite concluir que condición posee el usuario del" + `
    " " + `
    "dispositivo, es decir " + `
    "podremos lograr un diagnostico básico en base a los datos recopilados por el medidor de signos vitales"
$new3 = "La interpretación de los parámetros permite concluir que condición posee el usuario del dispositivo, es decir podremos lograr un diagnostico básico en base a los datos recopilados por el medidor de signos vitales"
Replace-Text $old3 $new3 | Out-Null

# Find the paragraph that now holds this merged text, then delete every
# empty paragraph after it through the end of the body (just before sectPr).
$targetParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -ge 0) {
        $len = $p.Range.End - $p.Range.Start
        if ($len -gt 150) {
            $targetParaIndex = $i
        }
    }
}

if ($targetParaIndex -gt 0 -and $targetParaIndex -lt $d.Paragraphs.Count) {
    $delStart = $d.Paragraphs.Item($targetParaIndex + 1).Range.Start
    $delEnd = $d.Paragraphs.Item($d.Paragraphs.Count).Range.End
    if ($delEnd -gt $delStart) {
        $delRange = $d.Range($delStart, $delEnd)
        $delRange.Delete()
    }
}

Write-Output "done"
